$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date column) for rows 2 to 111
# from Excel serial date 45202 (2023-10-03) to 45203 (2023-10-04),
# keeping the existing date formatting/style.
for ($row = 2; $row -le 111; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
